# Retrofit all the UI check cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Please Select" -> "Please Select One" values
$ws.Range("C2").Value = "Please Select One"
$ws.Range("D2").Value = "Please Select One"

# Trim the Language list down to English / French / Japanese only
$ws.Range("B3").Value = "French"
$ws.Range("B4").Value = "Japanese"

# Clear out the now-unused language rows (B5:B11), keep styling
$ws.Range("B5:B11").ClearContents()

# Update the active selection/cursor position to B14
$ws.Application.Goto($ws.Range("B14"))

$wb.Save()
